# Updates cryptos list with latest scraped prices/volumes (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.162.36"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'498.41"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "'0.0951"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "'4.73"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").Value = "2.671.34"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'22.63"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "54.141.62"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "2.271.12"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "'10.20"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'4.14"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'303.14"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'61.06"
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'0.149"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "'7.28"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").Value = "'170.81"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "0.0₃0691"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'17.69"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'0.941"
$ws.Range("E35").Value = "  +9.48%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'4.80"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'124.59"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "'238.62"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'0.371"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "'16.14"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "'4.65"
$ws.Range("E51").Value = "  -0.48%  "
